# Regenerate the localization-status report for archive:
#  - the handoff that used to be pending is now showing as "In Translation"
#    (was "Ready for handoff") on every sheet that surfaces the status, and
#  - the "Status"-ish columns (zh-cn / de-de on the Overview sheet, and the
#    Status column on each language sheet) shrink to fit the new, shorter
#    text.

$wb = $excel.ActiveWorkbook

# Swap the status text everywhere it appears (Overview E2/F2, and C2 on
# each language sheet).
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation") | Out-Null
}

# Shrink the now-narrower status columns to fit "In Translation".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.55   # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = 12.55   # column F (de-de)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.55        # column C (Status)

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.55        # column C (Status)
